$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor values updated
$ws.Range("B3").Value = 0.9930009045483686
$ws.Range("C3").Value = 0.9928417572466479
$ws.Range("D3").Value = 0.9932420033692001

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor with new values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9925446792382711
$ws.Range("C4").Value = 0.9931567092727579
$ws.Range("D4").Value = 0.9932328216354307

# Row 5: AdaBoostRegressor -> MLPRegressor with new values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9944932773543496
$ws.Range("C5").Value = 0.9948485036069429
$ws.Range("D5").Value = 0.9946160190017307
